$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133-234 down to 134-235
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new data record
$ws.Cells.Item(133, 1).Value = 8
$ws.Cells.Item(133, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(133, 3).Value = "Coquimbo"
$ws.Cells.Item(133, 4).Value = 45128
$ws.Cells.Item(133, 5).Value = 4
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100107
$ws.Cells.Item(133, 8).Value = "Otros"
$ws.Cells.Item(133, 9).Value = 100107002
$ws.Cells.Item(133, 10).Value = "Chirimoya"
$ws.Cells.Item(133, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 300
$ws.Cells.Item(133, 14).Value = 25000
$ws.Cells.Item(133, 15).Value = 26000
$ws.Cells.Item(133, 16).Value = 25500
$ws.Cells.Item(133, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(133, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(133, 19).Value = 2550
$ws.Cells.Item(133, 20).Value = 10
